$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title: "Iteration 2 Plan" -> "Iteration 3 Plan"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Iteration 2", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Iteration 3", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. High-level objectives bullet list (numId=1):
#    Rewrite the first bullet and drop the other two (plus the blank
#    paragraph that used to follow them), leaving a single bullet describing
#    the "most consumed consumable" feature.
# ---------------------------------------------------------------------------
$pObjective = $d.Paragraphs.Item(3)
$objRange = $d.Range($pObjective.Range.Start, $pObjective.Range.End - 1)
$objRange.Find.Execute("Create a chrome extension that consumes a consumable.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Create a chrome extension that displays the most consumed consumable", 2) | Out-Null

$pKeepStart = $d.Paragraphs.Item(4)
$pKeepEnd = $d.Paragraphs.Item(6)
$d.Range($pKeepStart.Range.Start, $pKeepEnd.Range.End).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3. "List of work items" list (numId=2):
#    Rewrite the first work item and delete the two that followed it
#    ("implements authentication." / "shares a consumable."), shifting the
#    remaining items up.
# ---------------------------------------------------------------------------
$pWorkItem = $d.Paragraphs.Item(5)
$workRange = $d.Range($pWorkItem.Range.Start, $pWorkItem.Range.End - 1)
$workRange.Find.Execute("Write code for chrome extension that consumes a consumable.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Write code for chrome extension that displays the most consumed consumable.", 2) | Out-Null

$pDelStart = $d.Paragraphs.Item(6)
$pDelEnd = $d.Paragraphs.Item(7)
$d.Range($pDelStart.Range.Start, $pDelEnd.Range.End).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 4. Move the hidden "_GoBack" bookmark: it used to sit right after
#    "Complete Use Case 6"; it now belongs on the blank paragraph that
#    follows "Prepare powerpoint and Demo for presentation" (paragraph 14
#    after the deletions above). Adding a bookmark with the same name
#    elsewhere relocates it (bookmark names are unique), so the old
#    occurrence disappears automatically.
# ---------------------------------------------------------------------------
$pBlank = $d.Paragraphs.Item(14)
$pBlank.Range.Bookmarks.Add("_GoBack") | Out-Null
